$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.962230384349823
$ws.Range("B1").Value = 0.7586323618888855
$ws.Range("C1").Value = 2.636927843093872
$ws.Range("D1").Value = 3.510997295379639
$ws.Range("E1").Value = 1.398565649986267
